# Generate Report for Handoff
# Updates the "8dd3e1c2-b665-4f78-97c6-03636ca610d5" row's handoff/handback
# timestamp on each sheet so that it gets its own freshly-generated value
# instead of re-using the timestamp from the row above it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column D = "Latest Handoff Date", row 6 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D6").Value = "2016-32-17 22:32:09"

# --- zh-cn sheet: column E = "Latest Handoff Datetime", row 6 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E6").Value = "2016-03-17 22:32:06"

# --- de-de sheet: column E = "Latest Handoff Datetime", row 6 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E6").Value = "2016-03-17 22:32:09"
